$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.740.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.545.11'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.91'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.61%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.38'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.88%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.765.14'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.543.86'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.69%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.511'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.705.60'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.19'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '212.51'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0689'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.33%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -5.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.99'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.65'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.49'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.87'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0460'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.71%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.335.85'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.30%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.53%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.19%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.32%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.523'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.75'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.799'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.70%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.43%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.58'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.73'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.71%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'mCoin'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.25'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.30%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.679.55'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.88'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₇0966'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.22%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.04%  '
